$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "2022-Q4" sheet, positioned right after "总计" ---
$anchor = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($anchor)
$q4.Name = "2022-Q4"

# Borrow header + column-A formatting from the neighboring quarter sheet
# (bold/centered/bordered s="2" style) so the new sheet matches its siblings.
$fmtSrc = $wb.Worksheets.Item("2022-Q3")
$fmtSrc.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$fmtSrc.Range("A2:A21").Copy()
$q4.Range("A2:A21").PasteSpecial(-4122)

# Header row
$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

# Columns B (fund code) and D:G (numeric-looking text) must stay TEXT,
# matching the inlineStr cells in the source data (keeps leading zeros
# like "007040" and trailing zeros like "4.60" intact).
$q4.Range("B2:B21").NumberFormat = "@"
$q4.Range("D2:G21").NumberFormat = "@"

# Fund-holdings data for 2022-Q4 (A=index, H=rank stay numeric)
$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(2,2).Value = "007040"
$q4.Cells.Item(2,3).Value = "新疆前海联合泳隆灵活配置混合C"
$q4.Cells.Item(2,4).Value = "7.66"
$q4.Cells.Item(2,5).Value = "93.78"
$q4.Cells.Item(2,6).Value = "3.45"
$q4.Cells.Item(2,7).Value = "0.2643"
$q4.Cells.Item(2,8).Value = 10
$q4.Cells.Item(3,1).Value = 1
$q4.Cells.Item(3,2).Value = "610004"
$q4.Cells.Item(3,3).Value = "信澳中小盘混合"
$q4.Cells.Item(3,4).Value = "4.60"
$q4.Cells.Item(3,5).Value = "91.92"
$q4.Cells.Item(3,6).Value = "4.68"
$q4.Cells.Item(3,7).Value = "0.2153"
$q4.Cells.Item(3,8).Value = 10
$q4.Cells.Item(4,1).Value = 2
$q4.Cells.Item(4,2).Value = "519013"
$q4.Cells.Item(4,3).Value = "海富通风格优势混合"
$q4.Cells.Item(4,4).Value = "3.23"
$q4.Cells.Item(4,5).Value = "92.73"
$q4.Cells.Item(4,6).Value = "2.29"
$q4.Cells.Item(4,7).Value = "0.0740"
$q4.Cells.Item(4,8).Value = 8
$q4.Cells.Item(5,1).Value = 3
$q4.Cells.Item(5,2).Value = "005933"
$q4.Cells.Item(5,3).Value = "新疆前海联合先进制造灵活配置混合A"
$q4.Cells.Item(5,4).Value = "0.90"
$q4.Cells.Item(5,5).Value = "90.73"
$q4.Cells.Item(5,6).Value = "5.44"
$q4.Cells.Item(5,7).Value = "0.0490"
$q4.Cells.Item(5,8).Value = 3
$q4.Cells.Item(6,1).Value = 4
$q4.Cells.Item(6,2).Value = "013051"
$q4.Cells.Item(6,3).Value = "汇泉臻心致远混合A"
$q4.Cells.Item(6,4).Value = "2.02"
$q4.Cells.Item(6,5).Value = "76.67"
$q4.Cells.Item(6,6).Value = "2.36"
$q4.Cells.Item(6,7).Value = "0.0477"
$q4.Cells.Item(6,8).Value = 8
$q4.Cells.Item(7,1).Value = 5
$q4.Cells.Item(7,2).Value = "013052"
$q4.Cells.Item(7,3).Value = "汇泉臻心致远混合C"
$q4.Cells.Item(7,4).Value = "1.37"
$q4.Cells.Item(7,5).Value = "76.67"
$q4.Cells.Item(7,6).Value = "2.36"
$q4.Cells.Item(7,7).Value = "0.0323"
$q4.Cells.Item(7,8).Value = 8
$q4.Cells.Item(8,1).Value = 6
$q4.Cells.Item(8,2).Value = "004128"
$q4.Cells.Item(8,3).Value = "新疆前海联合泳隆灵活配置混合A"
$q4.Cells.Item(8,4).Value = "0.84"
$q4.Cells.Item(8,5).Value = "93.78"
$q4.Cells.Item(8,6).Value = "3.45"
$q4.Cells.Item(8,7).Value = "0.0290"
$q4.Cells.Item(8,8).Value = 10
$q4.Cells.Item(9,1).Value = 7
$q4.Cells.Item(9,2).Value = "011523"
$q4.Cells.Item(9,3).Value = "新疆前海联合产业趋势混合A"
$q4.Cells.Item(9,4).Value = "0.60"
$q4.Cells.Item(9,5).Value = "83.30"
$q4.Cells.Item(9,6).Value = "4.44"
$q4.Cells.Item(9,7).Value = "0.0266"
$q4.Cells.Item(9,8).Value = 6
$q4.Cells.Item(10,1).Value = 8
$q4.Cells.Item(10,2).Value = "008180"
$q4.Cells.Item(10,3).Value = "同泰慧利混合A"
$q4.Cells.Item(10,4).Value = "0.48"
$q4.Cells.Item(10,5).Value = "93.91"
$q4.Cells.Item(10,6).Value = "5.53"
$q4.Cells.Item(10,7).Value = "0.0265"
$q4.Cells.Item(10,8).Value = 5
$q4.Cells.Item(11,1).Value = 9
$q4.Cells.Item(11,2).Value = "005212"
$q4.Cells.Item(11,3).Value = "汇安稳裕债券"
$q4.Cells.Item(11,4).Value = "2.14"
$q4.Cells.Item(11,5).Value = "22.17"
$q4.Cells.Item(11,6).Value = "0.87"
$q4.Cells.Item(11,7).Value = "0.0186"
$q4.Cells.Item(11,8).Value = 10
$q4.Cells.Item(12,1).Value = 10
$q4.Cells.Item(12,2).Value = "011524"
$q4.Cells.Item(12,3).Value = "新疆前海联合产业趋势混合C"
$q4.Cells.Item(12,4).Value = "0.36"
$q4.Cells.Item(12,5).Value = "83.30"
$q4.Cells.Item(12,6).Value = "4.44"
$q4.Cells.Item(12,7).Value = "0.0160"
$q4.Cells.Item(12,8).Value = 6
$q4.Cells.Item(13,1).Value = 11
$q4.Cells.Item(13,2).Value = "003238"
$q4.Cells.Item(13,3).Value = "新华外延增长主题灵活配置混合"
$q4.Cells.Item(13,4).Value = "0.51"
$q4.Cells.Item(13,5).Value = "85.17"
$q4.Cells.Item(13,6).Value = "2.51"
$q4.Cells.Item(13,7).Value = "0.0128"
$q4.Cells.Item(13,8).Value = 5
$q4.Cells.Item(14,1).Value = 12
$q4.Cells.Item(14,2).Value = "008181"
$q4.Cells.Item(14,3).Value = "同泰慧利混合C"
$q4.Cells.Item(14,4).Value = "0.19"
$q4.Cells.Item(14,5).Value = "93.91"
$q4.Cells.Item(14,6).Value = "5.53"
$q4.Cells.Item(14,7).Value = "0.0105"
$q4.Cells.Item(14,8).Value = 5
$q4.Cells.Item(15,1).Value = 13
$q4.Cells.Item(15,2).Value = "012480"
$q4.Cells.Item(15,3).Value = "汇安信泰稳健一年持有期混合C"
$q4.Cells.Item(15,4).Value = "1.01"
$q4.Cells.Item(15,5).Value = "27.85"
$q4.Cells.Item(15,6).Value = "1.01"
$q4.Cells.Item(15,7).Value = "0.0102"
$q4.Cells.Item(15,8).Value = 9
$q4.Cells.Item(16,1).Value = 14
$q4.Cells.Item(16,2).Value = "012479"
$q4.Cells.Item(16,3).Value = "汇安信泰稳健一年持有期混合A"
$q4.Cells.Item(16,4).Value = "0.84"
$q4.Cells.Item(16,5).Value = "27.85"
$q4.Cells.Item(16,6).Value = "1.01"
$q4.Cells.Item(16,7).Value = "0.0085"
$q4.Cells.Item(16,8).Value = 9
$q4.Cells.Item(17,1).Value = 15
$q4.Cells.Item(17,2).Value = "014014"
$q4.Cells.Item(17,3).Value = "招商臻选平衡混合A"
$q4.Cells.Item(17,4).Value = "0.25"
$q4.Cells.Item(17,5).Value = "66.99"
$q4.Cells.Item(17,6).Value = "2.43"
$q4.Cells.Item(17,7).Value = "0.0061"
$q4.Cells.Item(17,8).Value = 9
$q4.Cells.Item(18,1).Value = 16
$q4.Cells.Item(18,2).Value = "005934"
$q4.Cells.Item(18,3).Value = "新疆前海联合先进制造灵活配置混合C"
$q4.Cells.Item(18,4).Value = "0.09"
$q4.Cells.Item(18,5).Value = "90.73"
$q4.Cells.Item(18,6).Value = "5.44"
$q4.Cells.Item(18,7).Value = "0.0049"
$q4.Cells.Item(18,8).Value = 3
$q4.Cells.Item(19,1).Value = 17
$q4.Cells.Item(19,2).Value = "014015"
$q4.Cells.Item(19,3).Value = "招商臻选平衡混合C"
$q4.Cells.Item(19,4).Value = "0.19"
$q4.Cells.Item(19,5).Value = "66.99"
$q4.Cells.Item(19,6).Value = "2.43"
$q4.Cells.Item(19,7).Value = "0.0046"
$q4.Cells.Item(19,8).Value = 9
$q4.Cells.Item(20,1).Value = 18
$q4.Cells.Item(20,2).Value = "010487"
$q4.Cells.Item(20,3).Value = "中银顺盈回报一年持有期混合"
$q4.Cells.Item(20,4).Value = "0.75"
$q4.Cells.Item(20,5).Value = "21.31"
$q4.Cells.Item(20,6).Value = "0.53"
$q4.Cells.Item(20,7).Value = "0.0040"
$q4.Cells.Item(20,8).Value = 10
$q4.Cells.Item(21,1).Value = 19
$q4.Cells.Item(21,2).Value = "000822"
$q4.Cells.Item(21,3).Value = "东海美丽中国灵活配置混合"
$q4.Cells.Item(21,4).Value = "0.12"
$q4.Cells.Item(21,5).Value = "76.79"
$q4.Cells.Item(21,6).Value = "2.59"
$q4.Cells.Item(21,7).Value = "0.0031"
$q4.Cells.Item(21,8).Value = 6

# --- 2. Prepend a matching "2022-Q4" row into the "总计" summary sheet ---
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()
$total.Cells.Item(3,1).Copy()
$total.Cells.Item(2,1).PasteSpecial(-4122)
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 20
$total.Cells.Item(2,4).Value = 0.86

# Column A is the running 0-based row index (row N -> N-2); re-stamp it for
# every data row now that a row was inserted at the top.
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(6,1).Value = 4
$total.Cells.Item(7,1).Value = 5
$total.Cells.Item(8,1).Value = 6

$total.Select()
$total.Range("A1").Select()
